{"js": "// Replace the date line and each of the 25 \"NN\u00f7N=NN, N\" answer cells with\n// their updated values. Every old string is unique in the document at the\n// moment it is searched (verified against the source diff order), so a\n// plain exact-text search + whole-match replace is safe and unambiguous.\nconst replacements = [\n  [\"2024-03-30 Saturday\", \"2024-03-31 Sunday\"],\n  [\"25\u00f75=5, 0\", \"97\u00f78=12, 1\"],\n  [\"79\u00f77=11, 2\", \"10\u00f76=1, 4\"],\n  [\"20\u00f78=2, 4\", \"11\u00f74=2, 3\"],\n  [\"77\u00f76=12, 5\", \"98\u00f73=32, 2\"],\n  [\"42\u00f78=5, 2\", \"36\u00f77=5, 1\"],\n  [\"63\u00f72=31, 1\", \"25\u00f75=5, 0\"],\n  [\"84\u00f75=16, 4\", \"61\u00f77=8, 5\"],\n  [\"34\u00f74=8, 2\", \"21\u00f73=7, 0\"],\n  [\"80\u00f77=11, 3\", \"22\u00f75=4, 2\"],\n  [\"64\u00f77=9, 1\", \"25\u00f79=2, 7\"],\n  [\"77\u00f79=8, 5\", \"36\u00f77=5, 1\"],\n  [\"11\u00f72=5, 1\", \"19\u00f73=6, 1\"],\n  [\"92\u00f79=10, 2\", \"59\u00f75=11, 4\"],\n  [\"30\u00f73=10, 0\", \"89\u00f75=17, 4\"],\n  [\"39\u00f75=7, 4\", \"53\u00f79=5, 8\"],\n  [\"52\u00f72=26, 0\", \"36\u00f77=5, 1\"],\n  [\"45\u00f72=22, 1\", \"77\u00f73=25, 2\"],\n  [\"35\u00f74=8, 3\", \"97\u00f72=48, 1\"],\n  [\"92\u00f75=18, 2\", \"44\u00f78=5, 4\"],\n  [\"32\u00f75=6, 2\", \"15\u00f74=3, 3\"],\n  [\"46\u00f76=7, 4\", \"89\u00f73=29, 2\"],\n  [\"43\u00f72=21, 1\", \"35\u00f72=17, 1\"],\n  [\"70\u00f75=14, 0\", \"17\u00f79=1, 8\"],\n  [\"54\u00f72=27, 0\", \"50\u00f72=25, 0\"],\n  [\"65\u00f79=7, 2\", \"46\u00f73=15, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the 25 \"NN\u00f7N=NN, N\" answer cells with\n# their updated values. Every \"Old\" string is unique in the document at the\n# moment it is searched (verified against the source diff order), so a\n# Find/Replace scoped to the whole document body is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"2024-03-30 Saturday\"; New = \"2024-03-31 Sunday\" },\n    @{ Old = \"25\u00f75=5, 0\";  New = \"97\u00f78=12, 1\" },\n    @{ Old = \"79\u00f77=11, 2\"; New = \"10\u00f76=1, 4\" },\n    @{ Old = \"20\u00f78=2, 4\";  New = \"11\u00f74=2, 3\" },\n    @{ Old = \"77\u00f76=12, 5\"; New = \"98\u00f73=32, 2\" },\n    @{ Old = \"42\u00f78=5, 2\";  New = \"36\u00f77=5, 1\" },\n    @{ Old = \"63\u00f72=31, 1\"; New = \"25\u00f75=5, 0\" },\n    @{ Old = \"84\u00f75=16, 4\"; New = \"61\u00f77=8, 5\" },\n    @{ Old = \"34\u00f74=8, 2\";  New = \"21\u00f73=7, 0\" },\n    @{ Old = \"80\u00f77=11, 3\"; New = \"22\u00f75=4, 2\" },\n    @{ Old = \"64\u00f77=9, 1\";  New = \"25\u00f79=2, 7\" },\n    @{ Old = \"77\u00f79=8, 5\";  New = \"36\u00f77=5, 1\" },\n    @{ Old = \"11\u00f72=5, 1\";  New = \"19\u00f73=6, 1\" },\n    @{ Old = \"92\u00f79=10, 2\"; New = \"59\u00f75=11, 4\" },\n    @{ Old = \"30\u00f73=10, 0\"; New = \"89\u00f75=17, 4\" },\n    @{ Old = \"39\u00f75=7, 4\";  New = \"53\u00f79=5, 8\" },\n    @{ Old = \"52\u00f72=26, 0\"; New = \"36\u00f77=5, 1\" },\n    @{ Old = \"45\u00f72=22, 1\"; New = \"77\u00f73=25, 2\" },\n    @{ Old = \"35\u00f74=8, 3\";  New = \"97\u00f72=48, 1\" },\n    @{ Old = \"92\u00f75=18, 2\"; New = \"44\u00f78=5, 4\" },\n    @{ Old = \"32\u00f75=6, 2\";  New = \"15\u00f74=3, 3\" },\n    @{ Old = \"46\u00f76=7, 4\";  New = \"89\u00f73=29, 2\" },\n    @{ Old = \"43\u00f72=21, 1\"; New = \"35\u00f72=17, 1\" },\n    @{ Old = \"70\u00f75=14, 0\"; New = \"17\u00f79=1, 8\" },\n    @{ Old = \"54\u00f72=27, 0\"; New = \"50\u00f72=25, 0\" },\n    @{ Old = \"65\u00f79=7, 2\";  New = \"46\u00f73=15, 1\" }\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p.Old\n    $find.Replacement.Text = $p.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
